$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

Set-TextValue "D2" '30.109.23'
$ws.Range("E2").Value = '  -1.48%  '
Set-TextValue "D3" '2.108.06'
$ws.Range("E3").Value = '  -0.03%  '
Set-TextValue "D4" '1.003'
$ws.Range("E4").Value = '  -0.85%  '
Set-TextValue "D5" '350.13'
$ws.Range("E5").Value = '  +4.19%  '
Set-TextValue "D6" '1.003'
Set-TextValue "D7" '0.5164'
$ws.Range("E7").Value = '  -1.58%  '
$ws.Range("E8").Value = '  -1.29%  '
Set-TextValue "D9" '52.76'
$ws.Range("E9").Value = '  -5.51%  '
Set-TextValue "D10" '0.08965'
$ws.Range("E10").Value = '  -0.87%  '
$ws.Range("E11").Value = '  +0.78%  '
Set-TextValue "D12" '25.93'
$ws.Range("E12").Value = '  +5.56%  '
Set-TextValue "D13" '2.105.68'
$ws.Range("E13").Value = '  -0.67%  '
Set-TextValue "D14" '8.269'
$ws.Range("E14").Value = '  +2.01%  '
Set-TextValue "D15" '6.760'
$ws.Range("E15").Value = '  -1.23%  '
Set-TextValue "D16" '99.38'
$ws.Range("E16").Value = '  +2.25%  '
Set-TextValue "D17" '0.00001148'
$ws.Range("E17").Value = '  -2.80%  '
$ws.Range("E18").Value = '  -0.78%  '
Set-TextValue "D19" '20.90'
$ws.Range("E19").Value = '  +8.16%  '
Set-TextValue "D20" '0.06670'
$ws.Range("E20").Value = '  -0.30%  '
Set-TextValue "D21" '1.002'
$ws.Range("E21").Value = '  -0.75%  '
Set-TextValue "D22" '6.314'
$ws.Range("E22").Value = '  +0.91%  '
Set-TextValue "D23" '30.199.65'
$ws.Range("E23").Value = '  -1.39%  '
Set-TextValue "D24" '12.92'
Set-TextValue "D25" '2.353'
$ws.Range("E25").Value = '  -0.41%  '
Set-TextValue "D26" '2.357.09'
$ws.Range("E26").Value = '  -0.35%  '
Set-TextValue "D27" '22.10'
Set-TextValue "D28" '2.570'
$ws.Range("E28").Value = '  +2.22%  '
Set-TextValue "D29" '163.13'
$ws.Range("E29").Value = '  -0.15%  '
Set-TextValue "D30" '133.81'
$ws.Range("E30").Value = '  +0.22%  '
Set-TextValue "D31" '1.184'
$ws.Range("E31").Value = '  -2.64%  '
Set-TextValue "D32" '0.1069'
$ws.Range("E32").Value = '  +0.09%  '
Set-TextValue "D33" '1.651'
$ws.Range("E33").Value = '  +2.02%  '
Set-TextValue "D34" '6.281'
$ws.Range("E34").Value = '  -0.86%  '
Set-TextValue "D35" '3.985'
$ws.Range("E35").Value = '  +0.45%  '
Set-TextValue "D36" '5.917'
$ws.Range("E36").Value = '  +0.94%  '
$ws.Range("E37").Value = '  -2.05%  '
Set-TextValue "D38" '0.02591'
$ws.Range("E38").Value = '  -0.81%  '
Set-TextValue "D39" '0.06849'
$ws.Range("E39").Value = '  +0.50%  '
Set-TextValue "D40" '0.2325'
$ws.Range("E40").Value = '  +0.64%  '
Set-TextValue "D41" '12.59'
$ws.Range("E41").Value = '  +0.07%  '
Set-TextValue "D42" '0.6860'
$ws.Range("E42").Value = '  +0.39%  '
$ws.Range("E43").Value = '  -0.48%  '
Set-TextValue "D44" '14.34'
$ws.Range("E44").Value = '  +2.40%  '
Set-TextValue "D45" '0.6438'
$ws.Range("E45").Value = '  +0.05%  '
Set-TextValue "D46" '2.312'
$ws.Range("E46").Value = '  +0.40%  '
$ws.Range("E47").Value = '  +3.96%  '
Set-TextValue "D48" '3.669'
$ws.Range("E48").Value = '  -0.13%  '
Set-TextValue "D49" '84.09'
$ws.Range("E49").Value = '  +1.15%  '
$ws.Range("E50").Value = '  -2.02%  '
Set-TextValue "D51" '0.07241'
$ws.Range("E51").Value = '  +0.79%  '
